$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 22; this shifts existing rows 22-51 down to 23-52.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new record's values.
$ws.Cells.Item(22, 1).Value = 5
$ws.Cells.Item(22, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(22, 3).Value = "Maule"
$ws.Cells.Item(22, 4).Value = 44571
$ws.Cells.Item(22, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(22, 5).Value = 7
$ws.Cells.Item(22, 6).Value = "Fruta"
$ws.Cells.Item(22, 7).Value = 100101
$ws.Cells.Item(22, 8).Value = "Berries"
$ws.Cells.Item(22, 9).Value = 100101001
$ws.Cells.Item(22, 10).Value = "Arándano (blue)"
$ws.Cells.Item(22, 11).Value = "Sin especificar"
$ws.Cells.Item(22, 12).Value = "Segunda"
$ws.Cells.Item(22, 13).Value = 120
$ws.Cells.Item(22, 14).Value = 3200
$ws.Cells.Item(22, 15).Value = 3200
$ws.Cells.Item(22, 16).Value = 3200
$ws.Cells.Item(22, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(22, 18).Value = "Provincia de Linares"
$ws.Cells.Item(22, 19).Value = 1600
$ws.Cells.Item(22, 20).Value = 2
